$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 233, shifting rows 233:269 down to 234:270
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new data entry
$ws.Range("A233").Value = 10
$ws.Range("B233").Value = "Vega Modelo de Temuco"
$ws.Range("C233").Value = "La Araucanía"
$ws.Range("D233").Value = 45131
$ws.Range("E233").Value = 9
$ws.Range("F233").Value = "Fruta"
$ws.Range("G233").Value = 100104
$ws.Range("H233").Value = "Frutos de pepita"
$ws.Range("I233").Value = 100104001
$ws.Range("J233").Value = "Granada"
$ws.Range("K233").Value = "Wonderfull"
$ws.Range("L233").Value = "Primera"
$ws.Range("M233").Value = 300
$ws.Range("N233").Value = 16000
$ws.Range("O233").Value = 16000
$ws.Range("P233").Value = 16000
$ws.Range("Q233").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R233").Value = "Provincia de Limarí"
$ws.Range("S233").Value = 1600
$ws.Range("T233").Value = 10
